# Update the "dSF" column (F) values on Sheet1 with the newly re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = -1
    3  = 2
    4  = 13
    5  = 2
    7  = -2
    8  = -1
    9  = -3
    12 = -1
    14 = -2
    15 = 4
    17 = 4
    18 = 9
    20 = -6
    21 = 3
    22 = 10
    24 = -1
    25 = -4
    26 = -3
    27 = 3
    28 = -1
    29 = 6
    30 = -2
    31 = 3
    32 = 6
    33 = -2
    34 = -4
    35 = -4
    36 = 1
    37 = 4
    38 = -3
    39 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
